# edit.ps1 - apply the two changes captured by the commit:
#   1. Slide 5's table switches from the embedded "Light Style 2 - Accent 1"
#      table style to the built-in "No Style, No Grid" table style.
#   2. The deck's theme colour scheme (Integral / "Red Violet") is replaced
#      by the stock "Office Theme" colour scheme ("Office").

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$s   = $p.Slides.Item(5)
$sh  = $s.Shapes.Item(2)
$sh.Table.ApplyStyle("{AA7941A7-055B-4F32-A775-41FEBC5F7E31}")

# --- 2. Theme colour scheme -------------------------------------------------
# Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as COM RGB() longs (0x00BBGGRR).
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
